$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(41,42),
    @(49,50),
    @(62,64),
    @(139,140),
    @(166,167),
    @(168,169),
    @(170,171),
    @(179,180),
    @(202,204),
    @(210,211),
    @(230,231),
    @(233,234),
    @(285,286),
    @(297,298),
    @(307,308),
    @(359,360),
    @(373,374),
    @(425,426),
    @(436,437),
    @(447,448),
    @(471,472),
    @(498,499),
    @(500,501)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $cell1 = $ws.Cells.Item($r1, 4)
    $cell2 = $ws.Cells.Item($r2, 4)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
}
